$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the generic "Project manager" / "Team member N" header labels
# with the real contributors' names.
$ws.Range("C1").Value = "Zyrynyl Melendres"
$ws.Range("D1").Value = "William Smallwood"
$ws.Range("E1").Value = "Euan Sutherland"
$ws.Range("F1").Value = "Niall Swan"
$ws.Range("G1").Value = "Patryk `tAugusewicz"
$ws.Range("H1").Value = "Thomas Heaton"

# Fill in the first week's data: the date range and each contributor's
# share of the work for that week.
$ws.Range("B2").Value = "Jan 20-Jan 27"
$ws.Range("C2").Value = 0.24
$ws.Range("D2").Value = 0.19
$ws.Range("F2").Value = 0.19
$ws.Range("G2").Value = 0.19
$ws.Range("H2").Value = 0.19

$wb.Application.Calculate()

# Leave the selection where the author last left it.
$null = $ws.Range("E5").Select()
